$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '62.740.24'
Set-TextValue 'E2' '  -1.67%  '
Set-TextValue 'D3' '3.029.86'
Set-TextValue 'E3' '  -2.04%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '582.65'
Set-TextValue 'D6' '149.23'
Set-TextValue 'E6' '  -4.93%  '
Set-TextValue 'E8' '  -3.18%  '
Set-TextValue 'D9' '3.031.27'
Set-TextValue 'E9' '  -1.92%  '
Set-TextValue 'E10' '  -4.02%  '
Set-TextValue 'E11' '  -3.39%  '
Set-TextValue 'D12' '0.444'
Set-TextValue 'E12' '  -2.59%  '
Set-TextValue 'E13' '  -4.26%  '
Set-TextValue 'D14' '35.31'
Set-TextValue 'E14' '  -6.03%  '
Set-TextValue 'E15' '  +1.31%  '
Set-TextValue 'D16' '3.532.03'
Set-TextValue 'E16' '  -2.00%  '
Set-TextValue 'B17' 'Polkadot'
Set-TextValue 'C17' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D17' '7.05'
Set-TextValue 'E17' '  -1.86%  '
Set-TextValue 'B18' 'WrappedBTC'
Set-TextValue 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '62.715.70'
Set-TextValue 'E18' '  -1.63%  '
Set-TextValue 'D19' '3.028.57'
Set-TextValue 'E19' '  -2.02%  '
Set-TextValue 'D20' '468.42'
Set-TextValue 'E20' '  -2.57%  '
Set-TextValue 'D21' '14.06'
Set-TextValue 'E21' '  -4.04%  '
Set-TextValue 'D22' '0.692'
Set-TextValue 'E22' '  -3.01%  '
Set-TextValue 'D23' '7.39'
Set-TextValue 'E23' '  -2.54%  '
Set-TextValue 'D24' '2.38'
Set-TextValue 'E24' '  -2.31%  '
Set-TextValue 'D25' '81.09'
Set-TextValue 'E26' '  -3.75%  '
Set-TextValue 'D27' '10.45'
Set-TextValue 'E27' '  +1.31%  '
Set-TextValue 'E29' '  +0.00%  '
Set-TextValue 'D30' '7.22'
Set-TextValue 'E30' '  -4.14%  '
Set-TextValue 'E31' '  -2.33%  '
Set-TextValue 'E32' '  -1.55%  '
Set-TextValue 'D33' '27.54'
Set-TextValue 'E34' '  -5.51%  '
Set-TextValue 'E35' '  -1.28%  '
Set-TextValue 'D36' '0.0₃0799'
Set-TextValue 'E36' '  -6.80%  '
Set-TextValue 'D37' '5.78'
Set-TextValue 'E37' '  -5.05%  '
Set-TextValue 'E38' '  -3.36%  '
Set-TextValue 'D39' '50.25'
Set-TextValue 'E39' '  -1.35%  '
Set-TextValue 'D40' '2.97'
Set-TextValue 'E40' '  -15.59%  '
Set-TextValue 'D41' '9.01'
Set-TextValue 'E41' '  -3.90%  '
Set-TextValue 'D42' '422.07'
Set-TextValue 'D43' '0.282'
Set-TextValue 'E43' '  -2.47%  '
Set-TextValue 'D45' '2.794.04'
Set-TextValue 'E45' '  -1.39%  '
Set-TextValue 'D46' '0.0356'
Set-TextValue 'E46' '  -2.37%  '
Set-TextValue 'D47' '38.12'
Set-TextValue 'E47' '  -10.19%  '
Set-TextValue 'D48' '129.83'
Set-TextValue 'E48' '  -1.09%  '
Set-TextValue 'D50' '24.55'
Set-TextValue 'E50' '  -4.69%  '
Set-TextValue 'E51' '  -1.84%  '
